$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 311, shifting rows 311:340 down to 313:342.
$ws.Range("A311:A312").EntireRow.Insert()

# Populate new row 311 with a new weekly price entry.
$ws.Cells.Item(311, 1).Value = 6
$ws.Cells.Item(311, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(311, 3).Value = "Metropolitana"
$ws.Cells.Item(311, 4).Value = 44449
$ws.Cells.Item(311, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(311, 5).Value = 13
$ws.Cells.Item(311, 6).Value = 100112012
$ws.Cells.Item(311, 7).Value = "Espinaca"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 460
$ws.Cells.Item(311, 11).Value = 4500
$ws.Cells.Item(311, 12).Value = 5000
$ws.Cells.Item(311, 13).Value = 4728
$ws.Cells.Item(311, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(311, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(311, 16).Value = 473
$ws.Cells.Item(311, 17).Value = 10
$ws.Cells.Item(311, 18).Value = "Hortaliza"

# Populate new row 312 with a second new weekly price entry.
$ws.Cells.Item(312, 1).Value = 6
$ws.Cells.Item(312, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(312, 3).Value = "Metropolitana"
$ws.Cells.Item(312, 4).Value = 44449
$ws.Cells.Item(312, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(312, 5).Value = 13
$ws.Cells.Item(312, 6).Value = 100112012
$ws.Cells.Item(312, 7).Value = "Espinaca"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 430
$ws.Cells.Item(312, 11).Value = 4500
$ws.Cells.Item(312, 12).Value = 5000
$ws.Cells.Item(312, 13).Value = 4698
$ws.Cells.Item(312, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(312, 15).Value = "Región Metropolitana"
$ws.Cells.Item(312, 16).Value = 470
$ws.Cells.Item(312, 17).Value = 10
$ws.Cells.Item(312, 18).Value = "Hortaliza"
